$d = $word.ActiveDocument

# The final paragraph currently reads:
#   "Beva bu ... to be one of the most pupok opo of th the most popular thing on earth"
# (with a _GoBack bookmark right at the end, before </w:p>).
#
# It needs to become a paragraph split across four runs - inserting
# "absolely and most amazing" in place of the single word "most" (so the
# bookmark still sits between "and most amazing" and the trailing
# " pupok opo ..." text) - plus new paragraph-mark run properties
# (rtl / cs font hint) on the paragraph itself. Then a brand-new plain
# paragraph is appended after it.

$targetParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Beva bu because the Satan is the monkey problem is the most uglu u jug; ug; g ugly for tjr thr the mps most upstream to become yj the mere lady in the face pf pi outer universe the sca; scalable tory is the most poe prper poe powerflow thea that conclude into the d si sum p of the mpu mother hand to be one of the </w:t></w:r><w:r><w:t xml:space="preserve">absolely </w:t></w:r><w:r><w:t>and most amazing</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> pupok opo of th the most popular thing on earth</w:t></w:r></w:p>
'@

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>We are thinking of telling the truth' and the truth is to bexomcobex become the a spear head to the next level which e wo will apprea appear soon indeed</w:t></w:r></w:p>
'@

# Locate the paragraph by its distinctive lead-in text rather than a
# hard-coded index, so the script is resilient to minor structural
# differences.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Beva bu because the Satan")) {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs($targetIndex)
$p.Range.InsertXML($targetParaXml)

$p = $d.Paragraphs($targetIndex)
$p.Range.InsertParagraphAfter()

$newP = $d.Paragraphs($targetIndex + 1)
$newP.Range.InsertXML($newParaXml)
